# The deck currently has two theme parts:
#   ppt/theme/theme1.xml -> linked from the slide master ("Integral" palette)
#   ppt/theme/theme2.xml -> linked from the notes master   ("Office Theme" palette)
#
# The target edit swaps the two themes' content: the slide master's theme
# (theme1.xml, the one that actually paints every slide via schemeClr
# references) becomes the stock "Office Theme" palette, while the notes
# master keeps holding the old "Integral" palette bytes (theme2.xml).
#
# Font scheme and format scheme (fills / lines / effects) are identical
# between the two themes, so only the 12 colour-scheme slots - and the
# cosmetic name="" labels, where reachable - need to change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Slide.ThemeColorScheme is backed by the slide master's linked theme part
# (ppt/theme/theme1.xml). Index -> theme slot:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4
#   9 accent5  10 accent6  11 hlink  12 folHlink
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($index, $r, $g, $b) {
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (replaces the current Integral one)
Set-ThemeColor 1  0   0   0     # dk1      000000
Set-ThemeColor 2  255 255 255   # lt1      FFFFFF
Set-ThemeColor 3  68  84  106   # dk2      44546A
Set-ThemeColor 4  231 230 230   # lt2      E7E6E6
Set-ThemeColor 5  91  155 213   # accent1  5B9BD5
Set-ThemeColor 6  237 125 49    # accent2  ED7D31
Set-ThemeColor 7  165 165 165   # accent3  A5A5A5
Set-ThemeColor 8  255 192 0     # accent4  FFC000
Set-ThemeColor 9  68  114 196   # accent5  4472C4
Set-ThemeColor 10 112 173 71    # accent6  70AD47
Set-ThemeColor 11 5   99  193   # hlink    0563C1
Set-ThemeColor 12 149 79  114   # folHlink 954F72
